$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.641.80"
$ws.Range("D3").Value = "3.596.32"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'608.79"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'148.21"
$ws.Range("E6").Value = "  +2.40%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.136"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'8.05"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "4.205.96"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "'29.89"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "3.560.66"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "66.704.40"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "'11.49"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").Value = "'15.11"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "'427.55"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'0.619"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "'78.92"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "3.736.91"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").Value = "'8.29"
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'0.159"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").Value = "3.593.38"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "'25.46"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").Value = "'177.74"
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("D40").Value = "'0.0857"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "'5.22"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("E44").Value = "  +8.23%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D47").Value = "'25.05"
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("E48").Value = "  +3.87%  "
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "'0.952"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'0.235"
$ws.Range("E51").Value = "  -1.29%  "
